# entries-sample.xlsx: update the sample "Date*" entry in row 2 (A2) from
# 01/31/2019 to 09/01/2019 — a plain content edit to the sample data row,
# matching the upstream commit ("all assets file changed").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "09/01/2019"
